$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Split the opening address block into three short paragraphs:
#      "To,"  /  "Mitesh Mamtora"  /  "Pragati Land Developers" followed by
#      the existing line-break tail ("Kandivali West," / "Mumbai - 400 067"),
#      and give each of the three new paragraphs "space after" = 0.
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs(1).Range
$findText = "To Mitesh Bhai,`vPragati Land Developer"
$replText = "To,`rMitesh Mamtora`rPragati Land Developers"
$firstPara.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replText, 2) | Out-Null

$d.Paragraphs(1).SpaceAfter = 0
$d.Paragraphs(2).SpaceAfter = 0
$d.Paragraphs(3).SpaceAfter = 0

# ---------------------------------------------------------------------------
# 2) "Dear Mitesh Bhai," -> "Dear Mitesh,"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Mitesh Bhai", $true, $false, $false, $false, $false, $true, 1, $false, "Mitesh", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "Additionally, almost all the flats" -> "Additionally, all the flats"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("almost all", $true, $false, $false, $false, $false, $true, 1, $false, "all", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) "...inspect each flat and carry out the necessary..." ->
#    "...inspect each flat and conduct the necessary..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("carry out", $true, $false, $false, $false, $false, $true, 1, $false, "conduct", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Header: merge the "Hemukalani" / " Cross Road No. 2, ..." runs into a
#    single contiguous run (text content is unchanged - replace-with-self
#    forces Word to rebuild the run as one piece).
# ---------------------------------------------------------------------------
$hdr = $d.Sections(1).Headers(1)
$dash = [char]8211
$hdrLine = "Hemukalani Cross Road No. 2, Irani Wadi, Kandivali (West), Mumbai " + $dash + " 400067"
$hdr.Range.Find.Execute($hdrLine, $true, $false, $false, $false, $false, $true, 1, $false, $hdrLine, 2) | Out-Null
